$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowCount = 24

$colB = New-Object 'object[,]' $rowCount,1
$colB[0,0] = [double]0.9699704652908565
$colB[1,0] = [double]0.9313606958920388
$colB[2,0] = [double]0.9082488705025469
$colB[3,0] = [double]0.8989800109049497
$colB[4,0] = [double]0.8974499427688158
$colB[5,0] = [double]0.908123262717595
$colB[6,0] = [double]0.9565342725247774
$colB[7,0] = [double]1.056200432378915
$colB[8,0] = [double]1.132338090286964
$colB[9,0] = [double]1.167614299957506
$colB[10,0] = [double]1.181064986752062
$colB[11,0] = [double]1.178164029983748
$colB[12,0] = [double]1.168719043991189
$colB[13,0] = [double]1.162945750471863
$colB[14,0] = [double]1.130045619021075
$colB[15,0] = [double]1.11002672335627
$colB[16,0] = [double]1.098572687202534
$colB[17,0] = [double]1.094704898374374
$colB[18,0] = [double]1.112151526375669
$colB[19,0] = [double]1.171490757260926
$colB[20,0] = [double]1.210810771841579
$colB[21,0] = [double]1.189775607271258
$colB[22,0] = [double]1.111190731337842
$colB[23,0] = [double]1.028728601517713
$ws.Range("B2:B25").Value = $colB

$colC = New-Object 'object[,]' $rowCount,1
$colC[0,0] = [double]0.1266191015981235
$colC[1,0] = [double]0.120506148101839
$colC[2,0] = [double]0.1168355562707433
$colC[3,0] = [double]0.115360510851545
$colC[4,0] = [double]0.1151168320464819
$colC[5,0] = [double]0.1168155793617984
$colC[6,0] = [double]0.1244941302330034
$colC[7,0] = [double]0.1402126111896393
$colC[8,0] = [double]0.1521708087757929
$colC[9,0] = [double]0.1577014948059059
$colC[10,0] = [double]0.1598089868962802
$colC[11,0] = [double]0.1593545150404054
$colC[12,0] = [double]0.1578746156241948
$colC[13,0] = [double]0.1569698482380488
$colC[14,0] = [double]0.1518112000748602
$colC[15,0] = [double]0.1486698586175805
$colC[16,0] = [double]0.1468715877996374
$colC[17,0] = [double]0.1462641895681998
$colC[18,0] = [double]0.149003374766977
$colC[19,0] = [double]0.1583089406720148
$colC[20,0] = [double]0.1644673010187887
$colC[21,0] = [double]0.1611734298132887
$colC[22,0] = [double]0.1488525680502732
$colC[23,0] = [double]0.1358888823014013
$ws.Range("C2:C25").Value = $colC

$colD = New-Object 'object[,]' $rowCount,1
$colD[0,0] = [double]0.05724020510061223
$colD[1,0] = [double]0.05651891498867911
$colD[2,0] = [double]0.05611064452490666
$colD[3,0] = [double]0.05595300421003913
$colD[4,0] = [double]0.05592735655793746
$colD[5,0] = [double]0.05610848312905148
$colD[6,0] = [double]0.05698434055022972
$colD[7,0] = [double]0.0589752601796647
$colD[8,0] = [double]0.06060328924332481
$colD[9,0] = [double]0.06137955143382356
$colD[10,0] = [double]0.06167860344250897
$colD[11,0] = [double]0.06161397096561672
$colD[12,0] = [double]0.06140405261456294
$colD[13,0] = [double]0.06127613463623049
$colD[14,0] = [double]0.06055327375300124
$colD[15,0] = [double]0.06011893484584618
$colD[16,0] = [double]0.05987247319336575
$colD[17,0] = [double]0.05978960327738037
$colD[18,0] = [double]0.06016482357996011
$colD[19,0] = [double]0.0614655725916009
$colD[20,0] = [double]0.06234539076115908
$colD[21,0] = [double]0.06187310719842287
$colD[22,0] = [double]0.06014406717156362
$colD[23,0] = [double]0.05840753894162987
$ws.Range("D2:D25").Value = $colD

$colF = New-Object 'object[,]' $rowCount,1
$colF[0,0] = [double]3.371857794526562
$colF[1,0] = [double]3.351252924252506
$colF[2,0] = [double]3.339962804262981
$colF[3,0] = [double]3.335704069768113
$colF[4,0] = [double]3.335017564135399
$colF[5,0] = [double]3.339903984740729
$colF[6,0] = [double]3.364470524913756
$colF[7,0] = [double]3.423464322003198
$colF[8,0] = [double]3.473435298741549
$colF[9,0] = [double]3.497615423716297
$colF[10,0] = [double]3.50698051456385
$colF[11,0] = [double]3.504954291030174
$colF[12,0] = [double]3.498381712518352
$colF[13,0] = [double]3.494382995744672
$colF[14,0] = [double]3.471884241169136
$colF[15,0] = [double]3.458453131723388
$colF[16,0] = [double]3.450864189936155
$colF[17,0] = [double]3.448318099225077
$colF[18,0] = [double]3.459868789477895
$colF[19,0] = [double]3.500306573738271
$colF[20,0] = [double]3.527951076397585
$colF[21,0] = [double]3.513085289254377
$colF[22,0] = [double]3.459228357027257
$colF[23,0] = [double]3.406343351239116
$ws.Range("F2:F25").Value = $colF

$colG = New-Object 'object[,]' $rowCount,1
$colG[0,0] = [double]0.002581779450897264
$colG[1,0] = [double]0.00258606108607259
$colG[2,0] = [double]0.002588829293119789
$colG[3,0] = [double]0.002589992493431316
$colG[4,0] = [double]0.002590187767192664
$colG[5,0] = [double]0.002588844838039164
$colG[6,0] = [double]0.002583226922594273
$colG[7,0] = [double]0.002573310029240289
$colG[8,0] = [double]0.002566687291197121
$colG[9,0] = [double]0.002563816893465187
$colG[10,0] = [double]0.002562750295784809
$colG[11,0] = [double]0.002562979102975574
$colG[12,0] = [double]0.002563728736484708
$colG[13,0] = [double]0.002564190556970648
$colG[14,0] = [double]0.002566877733310176
$colG[15,0] = [double]0.002568562605021748
$colG[16,0] = [double]0.002569545100599134
$colG[17,0] = [double]0.002569880061798269
$colG[18,0] = [double]0.002568381861220326
$colG[19,0] = [double]0.002563507999107712
$colG[20,0] = [double]0.002560441272117272
$colG[21,0] = [double]0.002562067222699017
$colG[22,0] = [double]0.002568463532511103
$colG[23,0] = [double]0.002575875823025202
$ws.Range("G2:G25").Value = $colG

$colJ = New-Object 'object[,]' $rowCount,1
$colJ[0,0] = [double]0.2537046272738408
$colJ[1,0] = [double]0.2517523326817326
$colJ[2,0] = [double]0.2506700638824739
$colJ[3,0] = [double]0.2502583021302698
$colJ[4,0] = [double]0.2501916970584901
$colJ[5,0] = [double]0.2506643922128475
$colJ[6,0] = [double]0.2530072997707364
$colJ[7,0] = [double]0.2585267217695986
$colJ[8,0] = [double]0.2631479733677082
$colJ[9,0] = [double]0.265373776798242
$colJ[10,0] = [double]0.2662344302751052
$colJ[11,0] = [double]0.2660482816000922
$colJ[12,0] = [double]0.2654442266194366
$colJ[13,0] = [double]0.2650765429082043
$colJ[14,0] = [double]0.2630050002051831
$colJ[15,0] = [double]0.2617658400763503
$colJ[16,0] = [double]0.2610647379422062
$colJ[17,0] = [double]0.2608293536611157
$colJ[18,0] = [double]0.2618965469826264
$colJ[19,0] = [double]0.2656211692268045
$colJ[20,0] = [double]0.268159129129856
$colJ[21,0] = [double]0.2667950762828895
$colJ[22,0] = [double]0.2618374191715986
$colJ[23,0] = [double]0.2569343176167962
$ws.Range("J2:J25").Value = $colJ

$colK = New-Object 'object[,]' $rowCount,1
$colK[0,0] = [double]1.018363857192583
$colK[1,0] = [double]0.9761220786593867
$colK[2,0] = [double]0.9508190341340992
$colK[3,0] = [double]0.9406669266648464
$colK[4,0] = [double]0.9389907790089183
$colK[5,0] = [double]0.9506814754188611
$colK[6,0] = [double]1.003667236643111
$colK[7,0] = [double]1.112616739494683
$colK[8,0] = [double]1.19577157389233
$colK[9,0] = [double]1.234284215565452
$colK[10,0] = [double]1.248966923486762
$colK[11,0] = [double]1.245800340464342
$colK[12,0] = [double]1.235490188859387
$colK[13,0] = [double]1.229187799587635
$colK[14,0] = [double]1.19326849266119
$colK[15,0] = [double]1.171408833975107
$colK[16,0] = [double]1.158900202515781
$colK[17,0] = [double]1.154676058196088
$colK[18,0] = [double]1.173729157890108
$colK[19,0] = [double]1.238515846369552
$colK[20,0] = [double]1.281433752559025
$colK[21,0] = [double]1.258474850414103
$colK[22,0] = [double]1.172679956388805
$colK[23,0] = [double]1.082599436514357
$ws.Range("K2:K25").Value = $colK

$colM = New-Object 'object[,]' $rowCount,1
$colM[0,0] = [double]0.4264926444671744
$colM[1,0] = [double]0.4150881845968755
$colM[2,0] = [double]0.4083395416244215
$colM[3,0] = [double]0.4056531895677864
$colM[4,0] = [double]0.4052109738682219
$colM[5,0] = [double]0.4083030543582638
$colM[6,0] = [double]0.4225076988141439
$colM[7,0] = [double]0.4523798252027902
$colM[8,0] = [double]0.4755644205214935
$colM[9,0] = [double]0.4863823700458951
$colM[10,0] = [double]0.4905179375119104
$colM[11,0] = [double]0.4896255316784064
$colM[12,0] = [double]0.4867218228468744
$colM[13,0] = [double]0.484948303718312
$colM[14,0] = [double]0.4748629029273914
$colM[15,0] = [double]0.4687453232901362
$colM[16,0] = [double]0.4652521729519989
$colM[17,0] = [double]0.4640738337469443
$colM[18,0] = [double]0.4693939085298169
$colM[19,0] = [double]0.4875736524253114
$colM[20,0] = [double]0.4996827748366925
$colM[21,0] = [double]0.4931990641651751
$colM[22,0] = [double]0.4691006086167206
$colM[23,0] = [double]0.4440818859253639
$ws.Range("M2:M25").Value = $colM

$colN = New-Object 'object[,]' $rowCount,1
$colN[0,0] = [double]3.476407650065909
$colN[1,0] = [double]3.480996978762818
$colN[2,0] = [double]3.484575949773159
$colN[3,0] = [double]3.486225589350767
$colN[4,0] = [double]3.486511051581303
$colN[5,0] = [double]3.484597423603617
$colN[6,0] = [double]3.477831955762724
$colN[7,0] = [double]3.470614593291714
$colN[8,0] = [double]3.469016834917738
$colN[9,0] = [double]3.469098274099323
$colN[10,0] = [double]3.469245632201421
$colN[11,0] = [double]3.469208709096549
$colN[12,0] = [double]3.469108060444114
$colN[13,0] = [double]3.469061593151366
$colN[14,0] = [double]3.46902780184007
$colN[15,0] = [double]3.46921429406197
$colN[16,0] = [double]3.469397619732817
$colN[17,0] = [double]3.469472745282374
$colN[18,0] = [double]3.469186567649885
$colN[19,0] = [double]3.469134458759612
$colN[20,0] = [double]3.469779702292115
$colN[21,0] = [double]3.469373070280909
$colN[22,0] = [double]3.469198865697734
$colN[23,0] = [double]3.469061593151366
$ws.Range("N2:N25").Value = $colN

Write-Host "Updated pl_mw results block (columns B,C,D,F,G,J,K,M,N; rows 2-25)"